# Weekly fruit/vegetable price update: a new weekly price observation for
# "Espinaca" at "Vega Modelo de Temuco" needs to be inserted as row 148
# (shifting the existing rows 148:153 down to 149:154), matching the
# canonical diff (dimension grows from A1:R153 to A1:R154).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 148; Excel shifts 148:153 -> 149:154 and the
# new row inherits formatting from the row above (keeps date style "s=2"
# on column D).
$ws.Rows.Item(148).Insert()

# Populate the newly inserted row 148 with the new weekly record.
$ws.Range("A148").Value = 10
$ws.Range("B148").Value = "Vega Modelo de Temuco"
$ws.Range("C148").Value = "La Araucanía"
$ws.Range("D148").Value = 44747
$ws.Range("E148").Value = 9
$ws.Range("F148").Value = 100112012
$ws.Range("G148").Value = "Espinaca"
$ws.Range("H148").Value = "Sin especificar"
$ws.Range("I148").Value = "Primera"
$ws.Range("J148").Value = 35
$ws.Range("K148").Value = 10000
$ws.Range("L148").Value = 10000
$ws.Range("M148").Value = 10000
$ws.Range("N148").Value = "$/docena de atados"
$ws.Range("O148").Value = "Región de La Araucanía"
$ws.Range("P148").Value = 3333
$ws.Range("Q148").Value = 3
$ws.Range("R148").Value = "Hortaliza"
